# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the rows that just got (re)handed-off, and
# - Marks those same rows' Priority column as "ht" (handoff type) on the
#   per-language sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-19 00:22:37"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-19 00:22:31"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-19 00:22:37"
}
